$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Localización"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "Identificador"

# --- Row 2 (Juan) ---
$ws.Range("A2").Value = "Juan Torres Pardo"
$ws.Range("B2").Value = "22.971436; -43.182565"
$ws.Range("C2").Value = "juan@example.com"
$ws.Range("D2").Value = "87654321P"

# --- Row 3 (Luis) ---
$ws.Range("A3").Value = "Luis López Fernando"
$ws.Range("B3").Value = "32.97; -23.1"
$ws.Range("C3").Value = "luis@example.com"
$ws.Range("D3").Value = "19160962F"

# --- Row 4 (Ana) ---
$ws.Range("A4").Value = "Ana Torres Pardo"
$ws.Range("B4").Value = "21.26; 50,26"
$ws.Range("C4").Value = "ana@example.com"
$ws.Range("D4").Value = "09940449X"

# Columns E through G no longer hold data - clear their contents while
# leaving G4's formatting (it keeps its style, just loses its value).
$ws.Range("E1:G4").ClearContents()

# Selection now rests on D4 (last edited cell) instead of G1.
[void]$ws.Range("D4").Select()
